$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 corresponds to cm019. Update the syllabus entry:
# topic "Interactivity" -> "Building Shiny apps", and mark link_it TRUE.
$ws.Range("D20").Value = "Building Shiny apps"
$ws.Range("C20").Value = $true

# Update the active selection to C21 (matches the saved cursor position).
$ws.Range("C21").Select()
